# Apply the "LinuxForHealth" rebrand edit described by the commit diff.
#
# Summary of the change:
#   - Metadata sheet: URL / Version / Date / Publisher values updated
#     (ibm.com -> linuxforhealth.org, 7.0.0 -> 8.0.0, new publish date,
#     "Alvearie Team" -> "LinuxForHealth Team").
#   - Elements sheet: the Extension row's "Constraint(s)" cell (AI2),
#     which incorrectly duplicated the ele-1/ext-1 constraint text also
#     shown on the Extension.extension row (AI4), is cleared to blank.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/reference"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Clear the duplicated constraint text from the top-level Extension row;
# it remains (correctly) on the Extension.extension row (AI4).
$elements.Range("AI2").Value = ""
